$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 91410.83
$ws.Range("E2").Value = 35
$ws.Range("F2").Value = 7.99
$ws.Range("H2").Value = 56
